$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 5 data rows (27-31), which also trims the now-unused
# shared strings (crn_00025 .. crn_00029) and shrinks the dimension.
$ws.Rows("27:31").Delete()

# Update the measured values for rows 2-26 (columns B:I) with the new data.
$data = @{
    2 = @(1769.54131, 233.444182, 0.216583, 3.391931, 50.78125, 5, 130, 121)
    3 = @(1727.524844, 242.228542, 0.636426, 6.095082, 52.34375, 4, 134, 118)
    4 = @(1737.608153, 248.515757, 0.383505, 3.09084, 50, 4, 128, 124)
    5 = @(1725.150368, 243.014147, 0.609005, 5.060298, 46.875, 3, 120, 133)
    6 = @(1769.89245, 233.115228, 0.192106, 4.441251, 50, 3, 128, 125)
    7 = @(1741.383203, 246.587725, 0.079778, 3.868076, 52.34375, 5, 134, 117)
    8 = @(1729.469216, 261.227781, 0.334588, 4.782109, 49.21875, 4, 126, 126)
    9 = @(1757.89233, 245.706346, -0.024852, 4.386302, 51.5625, 5, 132, 119)
    10 = @(1744.902815, 227.481954, 0.296993, 3.917674, 51.5625, 4, 132, 120)
    11 = @(1752.741149, 236.532555, 0.324316, 4.643901, 53.90625, 4, 138, 114)
    12 = @(1722.217662, 253.147209, -0.715754, 1.596722, 54.6875, 4, 140, 112)
    13 = @(1794.144868, 232.026691, -0.206985, 2.523969, 51.5625, 4, 132, 120)
    14 = @(1745.685595, 267.662661, 0.315709, 4.220947, 50, 2, 128, 126)
    15 = @(1772.263861, 254.856945, 0.169575, 3.478333, 47.65625, 5, 122, 129)
    16 = @(1727.870403, 256.860981, 0.521527, 5.076763, 49.21875, 3, 126, 127)
    17 = @(1696.698131, 253.71489, 0.113652, 2.749505, 48.4375, 2, 124, 130)
    18 = @(1738.171894, 250.004131, -1.650949, 0.672688, 48.4375, 3, 124, 129)
    19 = @(1747.077781, 247.131886, -0.546278, 2.746638, 46.09375, 4, 118, 134)
    20 = @(1770.287538, 247.960486, 0.135606, 4.430412, 53.125, 5, 136, 115)
    21 = @(1797.91262, 219.138451, -1.034322, 1.791313, 49.21875, 4, 126, 126)
    22 = @(1740.775355, 248.524451, 0.083223, 4.635164, 49.21875, 4, 126, 126)
    23 = @(1739.08281, 258.080376, 0.005213, 4.272405, 47.65625, 3, 122, 131)
    24 = @(1743.294503, 242.917768, 0.572087, 6.088367, 51.5625, 3, 132, 121)
    25 = @(1754.704394, 241.550846, 0.460086, 4.040134, 51.5625, 3, 132, 121)
    26 = @(1763.622823, 237.730517, -0.250895, 2.487608, 49.21875, 3, 126, 127)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # Column B is index 2
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
